$d = $word.ActiveDocument

# Locate the "Docente(s) Responsável(eis)" heading paragraph and insert a
# new ListBullet paragraph right after it, naming the responsible instructor.
$found = $false
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Docente(s)*") {
        $p.Range.InsertParagraphAfter()
        $newPara = $p.Next()
        $newPara.Range.Text = "6712818 - Mauricio Lamano Ferreira"
        $newPara.Style = "ListBullet"
        $found = $true
        break
    }
}

Write-Output "inserted: $found"
